# Updated cryptos list on Sat Nov  4 05:46:36 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that sometimes look like plain numbers
# (e.g. "232.40", "41.09"). Excel's COM layer auto-coerces such strings
# into real numbers on assignment, which would both change the stored
# type and normalize away formatting (e.g. trailing zeros). Force the
# whole price column to Text first so every write below is kept as a
# literal string, then restore the default "Normal" style afterwards so
# the cells end up styled exactly as they started (no stray NumberFormat
# left behind on the data).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.940.80"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").Value = "1.840.45"
$ws.Range("E3").Value = "  +1.96%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "232.40"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("E6").Value = "  +2.51%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "41.09"
$ws.Range("E8").Value = "  +5.35%  "

$ws.Range("E9").Value = "  +3.25%  "

$ws.Range("D10").Value = "0.0691"
$ws.Range("E10").Value = "  +1.98%  "

$ws.Range("D11").Value = "0.0981"
$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").Value = "2.106.73"
$ws.Range("E12").Value = "  +1.96%  "

$ws.Range("D13").Value = "11.35"
$ws.Range("E13").Value = "  +4.51%  "

$ws.Range("D14").Value = "1.844.46"
$ws.Range("E14").Value = "  +2.03%  "

$ws.Range("E15").Value = "  +1.71%  "

$ws.Range("D16").Value = "4.66"
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("D17").Value = "34.926.10"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").Value = "69.82"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").Value = "0.0₃0789"
$ws.Range("E19").Value = "  +1.05%  "

$ws.Range("D20").Value = "239.92"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "12.15"
$ws.Range("E21").Value = "  +3.52%  "

$ws.Range("E22").Value = "  +2.20%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("E24").Value = "  +1.29%  "

$ws.Range("D25").Value = "171.92"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("E26").Value = "  +2.29%  "

$ws.Range("D27").Value = "17.41"
$ws.Range("E27").Value = "  +1.78%  "

$ws.Range("E28").Value = "  +3.87%  "

$ws.Range("D29").Value = "1.66"
$ws.Range("E29").Value = "  +9.84%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("E31").Value = "  +1.49%  "

$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("D33").Value = "3.91"
$ws.Range("E33").Value = "  -0.55%  "

$ws.Range("E34").Value = "  +22.32%  "

$ws.Range("D35").Value = "1.95"
$ws.Range("E35").Value = "  +10.85%  "

$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("D37").Value = "0.748"
$ws.Range("E37").Value = "  +8.59%  "

$ws.Range("E38").Value = "  +10.91%  "

$ws.Range("D39").Value = "89.76"
$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("E40").Value = "  +3.49%  "

$ws.Range("D41").Value = "1.337.93"
$ws.Range("E41").Value = "  +2.36%  "

$ws.Range("D42").Value = "14.61"
$ws.Range("E42").Value = "  +3.09%  "

$ws.Range("E43").Value = "  -2.11%  "

$ws.Range("E44").Value = "  +2.01%  "

$ws.Range("E45").Value = "  +3.53%  "

$ws.Range("D47").Value = "6.30"
$ws.Range("E47").Value = "  +2.97%  "

$ws.Range("D48").Value = "2.026.85"
$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("D49").Value = "10.98"
$ws.Range("E49").Value = "  +70.00%  "

$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "3.39"
$ws.Range("E51").Value = "  +15.90%  "

# Restore the default cell style across the price column so the text
# forcing above doesn't leave a visible NumberFormat/style change.
$ws.Range("D2:D51").Style = "Normal"
